$d = $word.ActiveDocument

# Replace the FILLER text with the first paragraph's real content.
$d.Content.Find.Execute("FILLER", $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "The Florida Polytechnic University SGA Department of External Affairs is comprised of the SGA Director of External Affairs and additional deputies, such as the Deputy of Communication and the Deputy of Governmental Affairs.",
                         2)

# Locate the (now updated) first paragraph so we can append the remaining
# paragraphs right after it, in order.
$p1 = $d.Paragraphs.First
$tail = $p1.Range
$tail.Collapse(0)

$lines = @(
    "The Department of External Affairs is responsible for overseeing external communications,",
    "overseeing University committee appointees, overseeing Governmental Legislative Affairs,",
    "Florida Student Association coordination, Capitol Improvement Trust Fund Designation",
    "Correspondence, and opportunities for students outside of the University.",
    "Contact: SGA- externalaffairs@floridapoly.edu"
)

foreach ($line in $lines) {
    $tail = $tail.Paragraphs.Last.Range
    $tail.Collapse(0)
    $newPara = $tail.InsertParagraphAfter()
    $tail = $d.Paragraphs.Last.Range
    $tail.Collapse(0)
    $tail.InsertBefore($line)
}
